$d = $word.ActiveDocument

# --- Step 1: cascading renames of function titles + descriptions -----------
# Each "Función …" list item (bold title run, ':' run, line-break, description
# run(s)) is updated to show the content that used to belong to the *next*
# item, except the 5th one, which takes on the content of "Función verificar
# recursos por ambiente" (the item that used to be two slots further down,
# since the items in between are being dropped). The replace is scoped to
# each paragraph's own Range so that renaming paragraph N never collides with
# matching text that a previous step may have just written into paragraph
# N-1.

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Paragraph 32: "Función calcular espacios disponibles" -> "Función contar dispositivos en minuta"
Replace-InParagraph 32 "Función calcular espacios disponibles" "Función contar dispositivos en minuta"
Replace-InParagraph 32 "Esta función devuelve la cantidad actual de espacios disponibles en el parqueadero, restando el número de vehículos ingresados a la capacidad total. Sirve para validar si hay cupo antes de autorizar un ingreso." "Esta función calcula cuántos dispositivos fueron entregados y cuántos fueron recibidos en una minuta específica. Su finalidad es detectar diferencias que puedan indicar pérdidas o robos."

# Paragraph 33: "Función contar dispositivos en minuta" -> "Función validar aprendiz activo"
Replace-InParagraph 33 "Función contar dispositivos en minuta" "Función validar aprendiz activo"
Replace-InParagraph 33 "Esta función calcula cuántos dispositivos fueron entregados y cuántos fueron recibidos en una minuta específica. Su finalidad es detectar diferencias que puedan indicar pérdidas o robos." "Esta función consulta si un aprendiz se encuentra activo en un programa de formación. Es útil antes de permitir el acceso o asignarle recursos dentro del sistema."

# Paragraph 34: "Función validar aprendiz activo" -> "Función calcular tiempo permanencia"
Replace-InParagraph 34 "Función validar aprendiz activo" "Función calcular tiempo permanencia"
Replace-InParagraph 34 "Esta función consulta si un aprendiz se encuentra activo en un programa de formación. Es útil antes de permitir el acceso o asignarle recursos dentro del sistema." "Esta función calcula el tiempo exacto que un estudiante, estuvo dentro del centro, usando la hora de entrada y salida. Esto apoya el control disciplinario y de seguridad."

# Paragraph 35: "Función calcular tiempo permanencia" -> "Función obtener estado ambiente"
Replace-InParagraph 35 "Función calcular tiempo permanencia" "Función obtener estado ambiente"
Replace-InParagraph 35 "Esta función calcula el tiempo exacto que un estudiante, visitante o funcionario estuvo dentro del centro, usando la hora de entrada y salida. Esto apoya el control disciplinario y de seguridad." "Esta función retorna el estado actual de un ambiente (disponible, en uso, mantenimiento), basándose en los horarios y reportes. Es clave para asignación de espacios."

# Paragraph 36: "Función obtener estado ambiente" -> "Función verificar recursos por ambiente"
Replace-InParagraph 36 "Función obtener estado ambiente" "Función verificar recursos por ambiente"
Replace-InParagraph 36 "Esta función retorna el estado actual de un ambiente (disponible, en uso, mantenimiento), basándose en los horarios y reportes. Es clave para asignación de espacios." "Esta función devuelve la cantidad y tipo de recursos que están asignados a un ambiente determinado, permitiendo validar si el inventario es correcto al generar minutas o asignaciones."

# --- Step 2: remove the now-superseded / not-required function items -------
# "Función listar visitantes del día", "Función verificar recursos por
# ambiente" (its original, further-down occurrence), "Función buscar
# responsable por visita" and "Función estado parqueadero" are removed
# entirely, leaving "Función obtener incidentes por ambiente" immediately
# after "Función verificar recursos por ambiente" (paragraph 36, renamed
# above).

$startRng = $d.Content
$startRng.Find.Execute("Función listar visitantes del día") | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("Función obtener incidentes por ambiente") | Out-Null
$endPos = $endRng.Start

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
